$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    "Q3" = 2.4
    "R3" = 1.53
    "G4" = 3.7
    "I4" = 2.25
    "J4" = 4.33
    "L4" = 3.1
    "M4" = 1.13
    "N4" = 6
    "U4" = 2.2
    "V4" = 1.62
    "AA4" = 34
    "AI4" = 9.5
    "AK4" = 21
    "AO4" = 21
    "AZ4" = 51
    "G5" = 1.27
    "I5" = 15
    "O5" = 1.36
    "P5" = 3
    "Q5" = 2.1
    "R5" = 1.7
    "U5" = 3.25
    "V5" = 1.33
    "W5" = 4.75
    "Y5" = 10
    "AD5" = 11
    "AE5" = 41
    "AH5" = 21
    "AJ5" = 41
    "AK5" = 251
    "AL5" = 126
    "AM5" = 151
    "AV5" = 126
    "AX5" = 67
    "G6" = 1.53
    "H6" = 3.9
    "I6" = 6.5
    "J6" = 2.1
    "M6" = 1.07
    "N6" = 9
    "Z6" = 10
    "AA6" = 13
    "AQ6" = 23
    "AS6" = 151
    "M8" = 1.05
    "N8" = 11
    "Q8" = 1.93
    "R8" = 1.93
    "G14" = 1.67
    "H14" = 3.7
    "I14" = 5.25
    "J14" = 2.3
    "L14" = 5.5
    "M14" = 1.07
    "N14" = 8.5
    "U14" = 2.1
    "V14" = 1.67
    "W14" = 6
    "X14" = 7
    "Z14" = 12
    "AB14" = 34
    "AC14" = 8.5
    "AE14" = 19
    "AF14" = 67
    "AI14" = 26
    "AJ14" = 17
    "AM14" = 51
    "AN14" = 3.5
    "AP14" = 23
    "AS14" = 201
    "AU14" = 9
    "AW14" = 7
    "AX14" = 29
    "AY14" = 41
    "BA14" = 151
    "BB14" = 351
}

foreach ($cell in $changes.Keys) {
    $ws.Range($cell).Value = $changes[$cell]
}
